# Macroferia Regional de Talca - Cebolla
# A new weekly price record was inserted at row 710 (pushing the existing
# rows 710-757 down to 711-758). Populate the new row with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 710; this shifts rows 710:757
# down to 711:758 (and extends the used range to row 758).
$ws.Rows(710).Insert()

# Fill in the newly inserted row 710 with the new record's values.
$ws.Range("A710").Value = 5
$ws.Range("B710").Value = "Macroferia Regional de Talca"
$ws.Range("C710").Value = "Maule"
$ws.Range("D710").Value = 44931
$ws.Range("E710").Value = 7
$ws.Range("F710").Value = 100112004
$ws.Range("G710").Value = "Cebolla"
$ws.Range("H710").Value = "Sin especificar"
$ws.Range("I710").Value = "1a nueva(o)"
$ws.Range("J710").Value = 50000
$ws.Range("K710").Value = 2200
$ws.Range("L710").Value = 2200
$ws.Range("M710").Value = 2200
$ws.Range("N710").Value = "`$/paquete 10 unidades (volumen en unidades)"
$ws.Range("O710").Value = "Región del Maule"
$ws.Range("P710").Value = 220
$ws.Range("Q710").Value = 10
$ws.Range("R710").Value = "Hortaliza"
